$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $cellRef, $val)
    # Leading apostrophe forces Excel to treat numeric-looking strings as
    # literal text instead of auto-converting them to numbers; resetting
    # the style afterwards avoids leaving a stray quote-prefix style behind.
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

$updates = @(
    @{ Cell = 'D2'; Value = '28.217.93' },
    @{ Cell = 'E2'; Value = '  -0.51%  ' },
    @{ Cell = 'D3'; Value = '1.832.02' },
    @{ Cell = 'E3'; Value = '  +1.32%  ' },
    @{ Cell = 'D4'; Value = '1.000' },
    @{ Cell = 'E4'; Value = '  -0.13%  ' },
    @{ Cell = 'D5'; Value = '309.85' },
    @{ Cell = 'E5'; Value = '  -1.07%  ' },
    @{ Cell = 'E6'; Value = '  -0.19%  ' },
    @{ Cell = 'D7'; Value = '0.4955' },
    @{ Cell = 'E7'; Value = '  -3.92%  ' },
    @{ Cell = 'D8'; Value = '0.1018' },
    @{ Cell = 'E8'; Value = '  +29.38%  ' },
    @{ Cell = 'D9'; Value = '0.3936' },
    @{ Cell = 'E9'; Value = '  -1.70%  ' },
    @{ Cell = 'D10'; Value = '1.111' },
    @{ Cell = 'E10'; Value = '  -0.10%  ' },
    @{ Cell = 'D11'; Value = '41.16' },
    @{ Cell = 'E11'; Value = '  +0.65%  ' },
    @{ Cell = 'D12'; Value = '6.425' },
    @{ Cell = 'E12'; Value = '  +1.10%  ' },
    @{ Cell = 'D13'; Value = '20.67' },
    @{ Cell = 'E13'; Value = '  +1.22%  ' },
    @{ Cell = 'D14'; Value = '1.000' },
    @{ Cell = 'E14'; Value = '  -0.14%  ' },
    @{ Cell = 'D15'; Value = '1.829.28' },
    @{ Cell = 'E15'; Value = '  +1.01%  ' },
    @{ Cell = 'D16'; Value = '7.349' },
    @{ Cell = 'E16'; Value = '  +0.45%  ' },
    @{ Cell = 'D17'; Value = '0.00001151' },
    @{ Cell = 'E17'; Value = '  +5.93%  ' },
    @{ Cell = 'D18'; Value = '92.94' },
    @{ Cell = 'E18'; Value = '  +0.32%  ' },
    @{ Cell = 'D19'; Value = '0.06650' },
    @{ Cell = 'E19'; Value = '  +1.15%  ' },
    @{ Cell = 'D20'; Value = '0.9990' },
    @{ Cell = 'E20'; Value = '  -0.19%  ' },
    @{ Cell = 'D21'; Value = '17.25' },
    @{ Cell = 'E21'; Value = '  -0.32%  ' },
    @{ Cell = 'E22'; Value = '  +0.14%  ' },
    @{ Cell = 'D23'; Value = '28.257.75' },
    @{ Cell = 'E23'; Value = '  -0.45%  ' },
    @{ Cell = 'D24'; Value = '11.29' },
    @{ Cell = 'E24'; Value = '  +1.46%  ' },
    @{ Cell = 'D25'; Value = '2.247' },
    @{ Cell = 'E25'; Value = '  +0.90%  ' },
    @{ Cell = 'B26'; Value = 'WrappedliquidstakedEther2.0' },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' },
    @{ Cell = 'D26'; Value = '2.040.59' },
    @{ Cell = 'E26'; Value = '  +0.93%  ' },
    @{ Cell = 'E27'; Value = '  +1.21%  ' },
    @{ Cell = 'B28'; Value = 'Monero' },
    @{ Cell = 'C28'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Cell = 'D28'; Value = '157.59' },
    @{ Cell = 'E28'; Value = '  -2.00%  ' },
    @{ Cell = 'D29'; Value = '2.428' },
    @{ Cell = 'E29'; Value = '  +0.82%  ' },
    @{ Cell = 'D30'; Value = '126.46' },
    @{ Cell = 'E30'; Value = '  -1.40%  ' },
    @{ Cell = 'D31'; Value = '0.1048' },
    @{ Cell = 'E31'; Value = '  -4.94%  ' },
    @{ Cell = 'D32'; Value = '1.039' },
    @{ Cell = 'E32'; Value = '  -2.47%  ' },
    @{ Cell = 'E33'; Value = '  +0.49%  ' },
    @{ Cell = 'D34'; Value = '3.591' },
    @{ Cell = 'E34'; Value = '  -2.16%  ' },
    @{ Cell = 'D35'; Value = '0.06759' },
    @{ Cell = 'E35'; Value = '  -6.24%  ' },
    @{ Cell = 'B36'; Value = 'VeChain' },
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D36'; Value = '0.02366' },
    @{ Cell = 'E36'; Value = '  +1.00%  ' },
    @{ Cell = 'B37'; Value = 'FraxShare' },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Cell = 'D37'; Value = '9.035' },
    @{ Cell = 'E37'; Value = '  -0.95%  ' },
    @{ Cell = 'D38'; Value = '0.2149' },
    @{ Cell = 'E38'; Value = '  -1.47%  ' },
    @{ Cell = 'B39'; Value = 'InternetComputer(DFINITY)' },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' },
    @{ Cell = 'D39'; Value = '4.986' },
    @{ Cell = 'E39'; Value = '  -1.31%  ' },
    @{ Cell = 'B40'; Value = 'Aptos' },
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' },
    @{ Cell = 'D40'; Value = '11.44' },
    @{ Cell = 'E40'; Value = '  -1.33%  ' },
    @{ Cell = 'D41'; Value = '0.6227' },
    @{ Cell = 'E41'; Value = '  +0.45%  ' },
    @{ Cell = 'D42'; Value = '1.174' },
    @{ Cell = 'E42'; Value = '  +1.67%  ' },
    @{ Cell = 'D43'; Value = '0.9990' },
    @{ Cell = 'E43'; Value = '  -0.20%  ' },
    @{ Cell = 'D44'; Value = '13.21' },
    @{ Cell = 'E44'; Value = '  -0.31%  ' },
    @{ Cell = 'D45'; Value = '0.5936' },
    @{ Cell = 'E45'; Value = '  -1.00%  ' },
    @{ Cell = 'E46'; Value = '  -1.36%  ' },
    @{ Cell = 'E47'; Value = '  -2.59%  ' },
    @{ Cell = 'D48'; Value = '123.86' },
    @{ Cell = 'E48'; Value = '  -1.32%  ' },
    @{ Cell = 'D49'; Value = '1.950' },
    @{ Cell = 'E49'; Value = '  +1.22%  ' },
    @{ Cell = 'E50'; Value = '  -2.91%  ' },
    @{ Cell = 'D51'; Value = '1.120' },
    @{ Cell = 'E51'; Value = '  +4.29%  ' }
)

foreach ($u in $updates) {
    Set-TextCell $ws $u.Cell $u.Value
}

Write-Output ('Updated ' + $updates.Count + ' cells')
